$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for all data rows (2-82)
# from 2023-09-06 (45175) to 2023-09-08 (45177)
for ($row = 2; $row -le 82; $row++) {
    $ws.Cells.Item($row, 3).Value = 45177
}
